$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.892.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.13%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.503.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.73%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'535.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.12%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'134.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.47%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.25%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +2.32%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.504.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.60%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.03%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -2.73%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.73%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -1.58%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.948.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.43%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'58.718.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.00%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'22.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.59%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.510.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.23%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E20").Value = "'  +1.54%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'321.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.94%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.13%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.38%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'65.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.93%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.37%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +1.82%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -1.47%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +1.12%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +0.88%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'171.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.57%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +1.62%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'6.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.08%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.80%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +0.14%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.41%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.26%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.06%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.04%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +3.42%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.829"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.26%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'36.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.47%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'3.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.36%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'274.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.33%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'131.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.65%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -1.79%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.589"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.54%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0937"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.74%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +2.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +2.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'16.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.25%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.747.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.15%  "
$ws.Range("E51").Style = "Normal"

